$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1797
$ws.Range("I40").Value = 1897.6666
$ws.Range("J40").Value = 1495
$ws.Range("K40").Value = 1897.6666
$ws.Range("L40").Value = 1495
$ws.Range("M40").Value = -1722.6666
$ws.Range("N40").Value = -1845
$ws.Range("H96").Value = 584.85
$ws.Range("I96").Value = 428.22223
$ws.Range("J96").Value = 1994.5
$ws.Range("K96").Value = 1284.66669
$ws.Range("L96").Value = 5983.5
$ws.Range("M96").Value = 88.33330999999998
$ws.Range("N96").Value = -8729.5
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("N103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("M103").ClearContents()
$ws.Range("H106").Value = 17095.309
$ws.Range("I106").Value = 2899.6667
$ws.Range("J106").Value = 36453
$ws.Range("K106").Value = 2899.6667
$ws.Range("L106").Value = 36453
$ws.Range("M106").Value = -2268.6667
$ws.Range("N106").Value = -37715
$ws.Range("H138").Value = 2108.2444
$ws.Range("I138").Value = 1520.3636
$ws.Range("J138").Value = 3724.9167
$ws.Range("K138").Value = 4561.0908
$ws.Range("L138").Value = 11174.7501
$ws.Range("M138").Value = 578.9092000000001
$ws.Range("N138").Value = -21454.7501

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5757.6665
$ws.Range("I63").Value = 5710.75
$ws.Range("J63").Value = 5851.5
$ws.Range("K63").Value = 5710.75
$ws.Range("L63").Value = 5851.5
$ws.Range("M63").Value = -5024.75
$ws.Range("N63").Value = -7223.5
$ws.Range("H66").Value = 5757.6665
$ws.Range("I66").Value = 5710.75
$ws.Range("J66").Value = 5851.5
$ws.Range("K66").Value = 28553.75
$ws.Range("L66").Value = 29257.5
$ws.Range("M66").Value = -25121.75
$ws.Range("N66").Value = -36121.5
$ws.Range("H80").Value = 31000
$ws.Range("I80").Value = 10100
$ws.Range("J80").Value = 33322.223
$ws.Range("K80").Value = 10100
$ws.Range("L80").Value = 33322.223
$ws.Range("M80").Value = -9102
$ws.Range("N80").Value = -35318.223
$ws.Range("H83").Value = 31000
$ws.Range("I83").Value = 10100
$ws.Range("J83").Value = 33322.223
$ws.Range("K83").Value = 30300
$ws.Range("L83").Value = 99966.66899999999
$ws.Range("M83").Value = -25308
$ws.Range("N83").Value = -109950.669
$ws.Range("H102").Value = 5149.8
$ws.Range("I102").Value = 5187.375
$ws.Range("J102").Value = 4999.5
$ws.Range("K102").Value = 5187.375
$ws.Range("L102").Value = 4999.5
$ws.Range("M102").Value = -3565.375
$ws.Range("N102").Value = -8243.5
$ws.Range("H132").Value = 7660.72
$ws.Range("I132").Value = 3813.25
$ws.Range("J132").Value = 100000
$ws.Range("K132").Value = 11439.75
$ws.Range("L132").Value = 300000
$ws.Range("M132").Value = -8909.75
$ws.Range("N132").Value = -305060

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2566.4443
$ws.Range("I99").Value = 2566.4443
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2566.4443
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1068.4443
$ws.Range("H105").Value = 3519.4
$ws.Range("I105").Value = 3979.0715
$ws.Range("J105").Value = 2446.8333
$ws.Range("K105").Value = 3979.0715
$ws.Range("L105").Value = 2446.8333
$ws.Range("M105").Value = -2232.0715
$ws.Range("N105").Value = -5940.8333
$ws.Range("H134").Value = 7875.2573
$ws.Range("I134").Value = 4791.2666
$ws.Range("J134").Value = 26379.2
$ws.Range("K134").Value = 14373.7998
$ws.Range("L134").Value = 79137.60000000001
$ws.Range("M134").Value = -11838.7998
$ws.Range("N134").Value = -84207.60000000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 249.1579
$ws.Range("I22").Value = 230
$ws.Range("J22").Value = 290.66666
$ws.Range("K22").Value = 230
$ws.Range("L22").Value = 290.66666
$ws.Range("M22").Value = 120
$ws.Range("N22").Value = -990.66666
$ws.Range("H31").Value = 1252.4517
$ws.Range("I31").Value = 953.5
$ws.Range("J31").Value = 1796
$ws.Range("K31").Value = 953.5
$ws.Range("L31").Value = 1796
$ws.Range("M31").Value = -658.5
$ws.Range("N31").Value = -2386
$ws.Range("H34").Value = 1252.4517
$ws.Range("I34").Value = 953.5
$ws.Range("J34").Value = 1796
$ws.Range("K34").Value = 953.5
$ws.Range("L34").Value = 1796
$ws.Range("M34").Value = -751.5
$ws.Range("N34").Value = -2200
$ws.Range("H58").Value = 3912.1333
$ws.Range("I58").Value = 2690.68
$ws.Range("J58").Value = 10019.4
$ws.Range("K58").Value = 2690.68
$ws.Range("L58").Value = 10019.4
$ws.Range("M58").Value = -2487.68
$ws.Range("N58").Value = -10425.4
$ws.Range("H94").Value = 922.65
$ws.Range("I94").Value = 536.8
$ws.Range("J94").Value = 1051.2667
$ws.Range("K94").Value = 536.8
$ws.Range("L94").Value = 1051.2667
$ws.Range("M94").Value = -85.79999999999995
$ws.Range("N94").Value = -1953.2667
$ws.Range("H99").Value = 11251.667
$ws.Range("I99").Value = 9778.6
$ws.Range("J99").Value = 11514.714
$ws.Range("K99").Value = 9778.6
$ws.Range("L99").Value = 11514.714
$ws.Range("M99").Value = -8280.6
$ws.Range("N99").Value = -14510.714
$ws.Range("H122").Value = 2196.3
$ws.Range("I122").Value = 1674
$ws.Range("J122").Value = 2326.875
$ws.Range("K122").Value = 5022
$ws.Range("L122").Value = 6980.625
$ws.Range("M122").Value = -2572
$ws.Range("N122").Value = -11880.625
$ws.Range("H126").Value = 11251.667
$ws.Range("I126").Value = 9778.6
$ws.Range("J126").Value = 11514.714
$ws.Range("K126").Value = 29335.8
$ws.Range("L126").Value = 34544.142
$ws.Range("M126").Value = -26865.8
$ws.Range("N126").Value = -39484.142
$ws.Range("H132").Value = 29312.715
$ws.Range("I132").Value = 23583.4
$ws.Range("J132").Value = 35280.75
$ws.Range("K132").Value = 70750.20000000001
$ws.Range("L132").Value = 105842.25
$ws.Range("M132").Value = -68220.20000000001
$ws.Range("N132").Value = -110902.25
$ws.Range("H134").Value = 12016.467
$ws.Range("I134").Value = 9249.923000000001
$ws.Range("J134").Value = 29999
$ws.Range("K134").Value = 27749.769
$ws.Range("L134").Value = 89997
$ws.Range("M134").Value = -25214.769
$ws.Range("N134").Value = -95067
$ws.Range("H136").Value = 3912.1333
$ws.Range("I136").Value = 2690.68
$ws.Range("J136").Value = 10019.4
$ws.Range("K136").Value = 8072.039999999999
$ws.Range("L136").Value = 30058.2
$ws.Range("M136").Value = -5522.039999999999
$ws.Range("N136").Value = -35158.2

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 694.53845
$ws.Range("I8").Value = 694.53845
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 2083.61535
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1944.61535
$ws.Range("H23").Value = 355.36365
$ws.Range("I23").Value = 76.75
$ws.Range("J23").Value = 514.5714
$ws.Range("K23").Value = 230.25
$ws.Range("L23").Value = 1543.7142
$ws.Range("M23").Value = 4.75
$ws.Range("N23").Value = -2013.7142
$ws.Range("H70").Value = 937.5
$ws.Range("I70").Value = 937.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2812.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2497.5
$ws.Range("H73").Value = 937.5
$ws.Range("I73").Value = 937.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2812.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1720.5
$ws.Range("H103").Value = 1251.7778
$ws.Range("I103").Value = 953.6
$ws.Range("J103").Value = 1624.5
$ws.Range("K103").Value = 2860.8
$ws.Range("L103").Value = 4873.5
$ws.Range("M103").Value = -1981.8
$ws.Range("N103").Value = -6631.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 37666.668
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 37666.668
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 37666.668
$ws.Range("N47").Value = -38802.668
$ws.Range("H70").Value = 6851.143
$ws.Range("I70").Value = 6105.2
$ws.Range("J70").Value = 8716
$ws.Range("K70").Value = 6105.2
$ws.Range("L70").Value = 8716
$ws.Range("M70").Value = -5835.2
$ws.Range("N70").Value = -9256
$ws.Range("H73").Value = 6851.143
$ws.Range("I73").Value = 6105.2
$ws.Range("J73").Value = 8716
$ws.Range("K73").Value = 6105.2
$ws.Range("L73").Value = 8716
$ws.Range("M73").Value = -5169.2
$ws.Range("N73").Value = -10588
$ws.Range("H126").Value = 2333
$ws.Range("I126").Value = 2499.5
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 7498.5
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -5028.5
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 4358.087
$ws.Range("I132").Value = 3951.8823
$ws.Range("J132").Value = 5509
$ws.Range("K132").Value = 11855.6469
$ws.Range("L132").Value = 16527
$ws.Range("M132").Value = -9325.6469
$ws.Range("N132").Value = -21587

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 79994.25
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 79994.25
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 79994.25
$ws.Range("N6").Value = -80218.25
$ws.Range("H46").Value = 1224.4706
$ws.Range("I46").Value = 1101.3334
$ws.Range("J46").Value = 1250.8572
$ws.Range("K46").Value = 1101.3334
$ws.Range("L46").Value = 1250.8572
$ws.Range("M46").Value = -913.3334
$ws.Range("N46").Value = -1626.8572
$ws.Range("H104").Value = 14590.714
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 14590.714
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 14590.714
$ws.Range("N104").Value = -21578.714

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("N68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("H70").Value = 36323.637
$ws.Range("I70").Value = 24746.5
$ws.Range("J70").Value = 38896.332
$ws.Range("K70").Value = 24746.5
$ws.Range("L70").Value = 38896.332
$ws.Range("M70").Value = -24431.5
$ws.Range("N70").Value = -39526.332
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("N71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("H73").Value = 36323.637
$ws.Range("I73").Value = 24746.5
$ws.Range("J73").Value = 38896.332
$ws.Range("K73").Value = 24746.5
$ws.Range("L73").Value = 38896.332
$ws.Range("M73").Value = -23654.5
$ws.Range("N73").Value = -41080.332
$ws.Range("H74").Value = 22331.334
$ws.Range("I74").Value = 29996
$ws.Range("J74").Value = 20798.4
$ws.Range("K74").Value = 29996
$ws.Range("L74").Value = 20798.4
$ws.Range("M74").Value = -29060
$ws.Range("N74").Value = -22670.4
$ws.Range("H77").Value = 22331.334
$ws.Range("I77").Value = 29996
$ws.Range("J77").Value = 20798.4
$ws.Range("K77").Value = 89988
$ws.Range("L77").Value = 62395.2
$ws.Range("M77").Value = -85308
$ws.Range("N77").Value = -71755.20000000001
$ws.Range("H100").Value = 568.8570999999999
$ws.Range("I100").Value = 568.8570999999999
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1137.7142
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -596.7141999999999
$ws.Range("H132").Value = 12621.725
$ws.Range("I132").Value = 6437.6187
$ws.Range("J132").Value = 33985
$ws.Range("K132").Value = 19312.8561
$ws.Range("L132").Value = 101955
$ws.Range("M132").Value = -16782.8561
$ws.Range("N132").Value = -107015
$ws.Range("H136").Value = 982.86664
$ws.Range("I136").Value = 900.52
$ws.Range("J136").Value = 1394.6
$ws.Range("K136").Value = 2701.56
$ws.Range("L136").Value = 4183.799999999999
$ws.Range("M136").Value = -151.5599999999999
$ws.Range("N136").Value = -9283.799999999999
